# Corrected data cleaning for pre/post/total fixation data
# - Updates several numeric metrics on rows 3-7 (recalculated after the fix)
# - Clears the "Unnamed: 0" header label from A1
# - Removes the bold/centered/bordered header formatting from row 1
# - Deletes the trailing blank row 10 (dimension shrinks from A1:AM10 to A1:AM9)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (Revisit count) ---
$ws.Range("J3").Value2 = 71
$ws.Range("L3").Value2 = 46
$ws.Range("O3").Value2 = 45
$ws.Range("S3").Value2 = 92
$ws.Range("Y3").Value2 = 29
$ws.Range("AC3").Value2 = 2
$ws.Range("AD3").Value2 = 77

# --- Row 4 (Fixation count) ---
$ws.Range("J4").Value2 = 441
$ws.Range("L4").Value2 = 89
$ws.Range("O4").Value2 = 76
$ws.Range("S4").Value2 = 343
$ws.Range("Y4").Value2 = 48
$ws.Range("AC4").Value2 = 5
$ws.Range("AD4").Value2 = 313

# --- Row 5 (Dwell time (ms)) ---
$ws.Range("J5").Value2 = 119021.4
$ws.Range("L5").Value2 = 26384.65
$ws.Range("O5").Value2 = 23451.07
$ws.Range("S5").Value2 = 96145.37
$ws.Range("Y5").Value2 = 14556.38
$ws.Range("AC5").Value2 = 1535.02
$ws.Range("AD5").Value2 = 88825.74

# --- Row 6 (Dwell time (%)) ---
$ws.Range("B6").Value2 = 0.74
$ws.Range("C6").Value2 = 1.48
$ws.Range("D6").Value2 = 0.97
$ws.Range("F6").Value2 = 1.54
$ws.Range("G6").Value2 = 4.95
$ws.Range("H6").Value2 = 5.42
$ws.Range("I6").Value2 = 7.54
$ws.Range("J6").Value2 = 43.04
$ws.Range("K6").Value2 = 15.56
$ws.Range("L6").Value2 = 9.54
$ws.Range("M6").Value2 = 6.72
$ws.Range("N6").Value2 = 7.85
$ws.Range("O6").Value2 = 8.48
$ws.Range("P6").Value2 = 0.63
$ws.Range("R6").Value2 = 6.57
$ws.Range("S6").Value2 = 34.77
$ws.Range("T6").Value2 = 2.01
$ws.Range("U6").Value2 = 2.8
$ws.Range("V6").Value2 = 0.47
$ws.Range("W6").Value2 = 2.05
$ws.Range("X6").Value2 = 3
$ws.Range("Y6").Value2 = 5.26
$ws.Range("Z6").Value2 = 2.9
$ws.Range("AA6").Value2 = 1.58
$ws.Range("AB6").Value2 = 0.09
$ws.Range("AC6").Value2 = 0.56
$ws.Range("AD6").Value2 = 32.12
$ws.Range("AE6").Value2 = 0.1
$ws.Range("AF6").Value2 = 0.71
$ws.Range("AG6").Value2 = 1.14
$ws.Range("AH6").Value2 = 0.33
$ws.Range("AJ6").Value2 = 3.67
$ws.Range("AK6").Value2 = 0.37
$ws.Range("AL6").Value2 = 1.11
$ws.Range("AM6").Value2 = 2.29

# --- Row 7 (Fixation duration (ms)) ---
$ws.Range("J7").Value2 = 269.89
$ws.Range("L7").Value2 = 296.46
$ws.Range("O7").Value2 = 308.57
$ws.Range("S7").Value2 = 280.31
$ws.Range("Y7").Value2 = 303.26
$ws.Range("AC7").Value2 = 307
$ws.Range("AD7").Value2 = 283.79

# --- Header row cleanup: drop "Unnamed: 0" label and the bold/border style ---
$ws.Range("A1").ClearContents()
$ws.Range("A1:AM1").ClearFormats()

# --- Remove trailing blank row 10 entirely (shrinks used range to A1:AM9) ---
$ws.Rows.Item(10).Delete()
